# Insert a new weekly price record for "Femacal de La Calera - Espinaca" as
# row 276 (Fecha 2022-08-11 / serial 44784), shifting all subsequent rows
# down by one (old row 276 becomes 277, ..., old row 385 becomes 386).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 276 - this shifts rows 276:385 down to 277:386
# and extends the sheet dimension to A1:R386.
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new record's data.
$ws.Cells.Item(276, 1).Value = 3
$ws.Cells.Item(276, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(276, 3).Value = "Coquimbo"
$ws.Cells.Item(276, 4).Value = 44784
$ws.Cells.Item(276, 5).Value = 5
$ws.Cells.Item(276, 6).Value = 100112012
$ws.Cells.Item(276, 7).Value = "Espinaca"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 220
$ws.Cells.Item(276, 11).Value = 4500
$ws.Cells.Item(276, 12).Value = 5000
$ws.Cells.Item(276, 13).Value = 4750
$ws.Cells.Item(276, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(276, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(276, 16).Value = 1583
$ws.Cells.Item(276, 17).Value = 3
$ws.Cells.Item(276, 18).Value = "Hortaliza"
